$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide 1 / "Textfeld 4": fix the "einlogen" -> "einloggen" typo and
#    widen the textbox (it autosizes to fit the now-longer first line).
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$txtShape = $s1.Shapes.Item(2)
$tr = $txtShape.TextFrame.TextRange

# Original run 1 = "Wer darf sich in die Software " (30 chars)
# Original run 2 = "einlogen" (8 chars)
# Target: "Wer darf sich in die " | "Software " | "einloggen"
$run1 = $tr.Characters(1, 21)
$run1.Text = "Wer darf sich in die "

$run2 = $tr.Characters(22, 9)
$run2.Text = "Software "

$run3 = $tr.Characters(31, 8)
$run3.Text = "einloggen"

# Widen the textbox to fit the corrected / re-flowed text.
$txtShape.Width = 294.286062992126

# ---------------------------------------------------------------------
# 2) Slide master + every slide layout: the cached preview text of the
#    slide-number field changes from the old "<Nr.>" token to "<#>".
# ---------------------------------------------------------------------
$numberToken = [char]0x2039 + "#" + [char]0x203A

function Update-SlideNumberField($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $t = $shp.TextFrame.TextRange.Text
            if ($t -like "*Nr.*") {
                $shp.TextFrame.TextRange.Text = $numberToken
            }
        }
    }
}

$master = $p.SlideMaster
Update-SlideNumberField $master.Shapes

for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Update-SlideNumberField $layout.Shapes
}
